$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 407, pushing existing rows 407..508 down to 408..509.
$ws.Rows.Item(407).Insert()

# The new row 407 carries the same descriptive/categorical values as the row
# that is now at 408 (the former row 407), but with new measurement data.
$ws.Range("A407").Value = $ws.Range("A408").Value2
$ws.Range("B407").Value = $ws.Range("B408").Value2
$ws.Range("C407").Value = $ws.Range("C408").Value2
$ws.Range("D407").Value = 44551
$ws.Range("E407").Value = $ws.Range("E408").Value2
$ws.Range("F407").Value = $ws.Range("F408").Value2
$ws.Range("G407").Value = $ws.Range("G408").Value2
$ws.Range("H407").Value = $ws.Range("H408").Value2
$ws.Range("I407").Value = $ws.Range("I408").Value2
$ws.Range("J407").Value = 1500
$ws.Range("K407").Value = 16500
$ws.Range("L407").Value = 17000
$ws.Range("M407").Value = 16700
$ws.Range("N407").Value = $ws.Range("N408").Value2
$ws.Range("O407").Value = $ws.Range("O408").Value2
$ws.Range("P407").Value = 1670
$ws.Range("Q407").Value = $ws.Range("Q408").Value2
$ws.Range("R407").Value = $ws.Range("R408").Value2
